# "data updated & ph added": the NFT500MC50 label used ":" as a separator;
# it is renamed to use "_" instead (new shared-string "placeholder" entry,
# the old one is dropped once nothing else references it).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A8").Value = "NFT500MC50_25_25"

# Re-fit column A now that its longest label changed, same as Excel
# recalculating the bestFit column width after the content edit.
$ws.Columns.Item(1).AutoFit() | Out-Null

# Leave the saved cursor position where the author last clicked.
$ws.Range("F14").Select() | Out-Null
